$d = $word.ActiveDocument

# ---- Page margins: 1417/1701 twips -> 720 twips (36pt) on all sides ----
$ps = $d.PageSetup
$ps.TopMargin = 36
$ps.RightMargin = 36
$ps.BottomMargin = 36
$ps.LeftMargin = 36

# ---- Section 4 (Deletar Musica): paragraphs 24-26
#   'Entrada:' + 'ID da musica: 2 ...' + 'Saida Esperada ... removida.'
#   -> 'Entrada: Clicar sobre o Btn deletar...' + empty paragraph + 'Saida Esperada ... removida.' (no lastRenderedPageBreak)
# NB: the trailing paragraph must not be empty or InsertXML silently drops a paragraph break,
# so the unchanged 'Saida Esperada' text is folded into this same call.
$xmlEntradaDeletarEmptySaida = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Entrada:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Clicar sobre o </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Btn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> deletar musica presente na frente da musica que vai ser deletada e exibido na pagina  Listar Musica para este teste irei deletar a musica  “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Shape</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>of</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>You</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Saída Esperada: A música é removida do sistema e a lista de músicas é atualizada, sem a música removida.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$pStart = $d.Paragraphs.Item(24)
$pEnd = $d.Paragraphs.Item(26)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.InsertXML($xmlEntradaDeletarEmptySaida)

# ---- Section 3 (Editar Musica): paragraphs 20-21
#   'Novo nome: ...' + 'Saida Esperada ... nome ...'
#   -> 'Ano de Lancamento: 1956' + 'Saida Esperada ... Ano de Lancamento ...' ----
$xmlAno1956Saida = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Ano de Lançamento: </w:t></w:r><w:r><w:t>1956</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Saída Esperada: A música é atualizada com o novo </w:t></w:r><w:r><w:t>Ano de Lançamento</w:t></w:r><w:r><w:t xml:space="preserve"> e a lista de músicas é exibida com as informações atualizadas.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$pStart = $d.Paragraphs.Item(20)
$pEnd = $d.Paragraphs.Item(21)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.InsertXML($xmlAno1956Saida)

# ---- Section 3: paragraphs 18-19 'Entrada:' + 'ID da musica: 1 ...' -> merged 'Entrada: Clicar sobre o Btn editar ...' ----
$xmlEntradaEditar = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Entrada:</w:t></w:r><w:r><w:t xml:space="preserve"> Clicar sobre o </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Btn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> editar musica presente na frente da musica que vai ser editada e exibido na pagina  Listar Musica para este teste irei editar o ano da musica “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Day-O</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>The</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Banana </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Boat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Song</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">)” para o ano correto do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lancamento</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$pStart = $d.Paragraphs.Item(18)
$pEnd = $d.Paragraphs.Item(19)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.InsertXML($xmlEntradaEditar)

# ---- Section 2 (Listar Musicas): paragraph 13 Premissa -> add ' usando o BTN "Listar Musicas".' ----
$xmlPremissaListar = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Premissa: O usuário deseja visualizar a lista de músicas cadastradas</w:t></w:r><w:r><w:t xml:space="preserve"> usando o BTN “Listar </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Musicas</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p = $d.Paragraphs.Item(13)
$r = $p.Range
$r.InsertXML($xmlPremissaListar)

# ---- Section 1 (Cadastrar Musica): paragraph 10 'Ano de Lancamento: 1975' -> '... 1955' ----
$xmlAno1955 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Ano de Lançamento: </w:t></w:r><w:r><w:t>1955</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p = $d.Paragraphs.Item(10)
$r = $p.Range
$r.InsertXML($xmlAno1955)

# ---- paragraph 9 'Artista: "Queen"' -> 'Artista: "Harry Belafonte"' ----
$xmlArtista = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Artista: "</w:t></w:r><w:r><w:t>Harry Belafonte</w:t></w:r><w:r><w:t>"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p = $d.Paragraphs.Item(9)
$r = $p.Range
$r.InsertXML($xmlArtista)

# ---- paragraph 8 'Genero: "Rock"' -> 'Genero: "Calypso"' ----
$xmlGenero = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Gênero: "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Calypso</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p = $d.Paragraphs.Item(8)
$r = $p.Range
$r.InsertXML($xmlGenero)

# ---- paragraph 7 'Nome: "Bohemian Rhapsody"' -> 'Nome: "Day-O (The Banana Boat Song)"' ----
$xmlNome = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Nome: </w:t></w:r><w:r><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Day-O</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>The</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Banana </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Boat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Song</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p = $d.Paragraphs.Item(7)
$r = $p.Range
$r.InsertXML($xmlNome)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
